$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row containing account "005696595" / "CLUBE" / 17755.92
# This is row 5 (row 1 = header "Conta"/"Nome"/"Saldo")
$ws.Rows.Item(5).Delete()
